$p = $ppt.ActivePresentation

# --- Slide 15: "Inheritance.XXX" -> "InheritanceType.XXX" ---------------
$s15 = $p.Slides.Item(15)
$sh15 = $s15.Shapes.Item(1)
$tr15 = $sh15.TextFrame.TextRange

$old1 = "Inheritance.TABLE_PER_CLASS"
$new1 = "InheritanceType.TABLE_PER_CLASS"
$full15 = $tr15.Text
$idx1 = $full15.IndexOf($old1)
$rng1 = $tr15.Characters($idx1 + 1, $old1.Length)
$rng1.Text = $new1

$old2 = "Inheritance.JOINED"
$new2 = "InheritanceType.JOINED"
$full15b = $tr15.Text
$idx2 = $full15b.IndexOf($old2)
$rng2 = $tr15.Characters($idx2 + 1, $old2.Length)
$rng2.Text = $new2

# --- Slide 8: "Mobile" class field typo "Mobile mobile;" -> "Employee employee;" ---
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(7)
$tr8 = $sh8.TextFrame.TextRange

$full8 = $tr8.Text
$oldLead = "    Mobile "
$idx3 = $full8.LastIndexOf("    Mobile mobile;")
$rng3 = $tr8.Characters($idx3 + 1, $oldLead.Length)
$rng3.Text = "    Employee "

$full8b = $tr8.Text
$idx4 = $full8b.IndexOf("Employee mobile") + "Employee ".Length
$rng4 = $tr8.Characters($idx4 + 1, "mobile".Length)
$rng4.Text = "employee"
